$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date number format (from column D, style s="2") so the new row 459
# gets the same date-time style as the rest of column D.
$dateFormat = $ws.Range("D424").NumberFormat

# Update rows 424-458 (shift each row's data down by one, i.e. row N receives what
# used to be row N-1's data) and insert new row 459 (prior row 458's data), while
# row 424 receives a brand-new weekly record.

# Row 424
$ws.Range("D424").Value = 45223
$ws.Range("L424").Value = 'Primera'
$ws.Range("M424").Value = 100
$ws.Range("N424").Value = 12000
$ws.Range("O424").Value = 12000
$ws.Range("P424").Value = 12000
$ws.Range("R424").Value = 'Brasil'
$ws.Range("S424").Value = 3000

# Row 425
$ws.Range("D425").Value = 44487
$ws.Range("L425").Value = 'Primera'
$ws.Range("M425").Value = 80
$ws.Range("N425").Value = 7500
$ws.Range("O425").Value = 8000
$ws.Range("P425").Value = 7750
$ws.Range("R425").Value = 'Perú'
$ws.Range("S425").Value = 1938

# Row 426
$ws.Range("D426").Value = 44714
$ws.Range("L426").Value = 'Primera'
$ws.Range("M426").Value = 120
$ws.Range("N426").Value = 7500
$ws.Range("O426").Value = 8000
$ws.Range("P426").Value = 7750
$ws.Range("R426").Value = 'Perú'
$ws.Range("S426").Value = 1938

# Row 427
$ws.Range("D427").Value = 44778
$ws.Range("L427").Value = 'Primera'
$ws.Range("M427").Value = 120
$ws.Range("N427").Value = 13000
$ws.Range("O427").Value = 14000
$ws.Range("P427").Value = 13500
$ws.Range("R427").Value = 'Brasil'
$ws.Range("S427").Value = 3375

# Row 428
$ws.Range("D428").Value = 44754
$ws.Range("L428").Value = 'Primera'
$ws.Range("M428").Value = 120
$ws.Range("N428").Value = 9000
$ws.Range("O428").Value = 9000
$ws.Range("P428").Value = 9000
$ws.Range("R428").Value = 'Perú'
$ws.Range("S428").Value = 2250

# Row 429
$ws.Range("D429").Value = 45215
$ws.Range("L429").Value = 'Primera'
$ws.Range("M429").Value = 100
$ws.Range("N429").Value = 13000
$ws.Range("O429").Value = 13000
$ws.Range("P429").Value = 13000
$ws.Range("R429").Value = 'Brasil'
$ws.Range("S429").Value = 3250

# Row 430
$ws.Range("D430").Value = 44924
$ws.Range("L430").Value = 'Primera'
$ws.Range("M430").Value = 120
$ws.Range("N430").Value = 8000
$ws.Range("O430").Value = 8500
$ws.Range("P430").Value = 8250
$ws.Range("R430").Value = 'Brasil'
$ws.Range("S430").Value = 2062

# Row 431
$ws.Range("D431").Value = 44481
$ws.Range("L431").Value = 'Primera'
$ws.Range("M431").Value = 200
$ws.Range("N431").Value = 7500
$ws.Range("O431").Value = 8000
$ws.Range("P431").Value = 7750
$ws.Range("R431").Value = 'Perú'
$ws.Range("S431").Value = 1938

# Row 432
$ws.Range("D432").Value = 44999
$ws.Range("L432").Value = 'Primera'
$ws.Range("M432").Value = 120
$ws.Range("N432").Value = 8000
$ws.Range("O432").Value = 8500
$ws.Range("P432").Value = 8250
$ws.Range("R432").Value = 'Perú'
$ws.Range("S432").Value = 2062

# Row 433
$ws.Range("D433").Value = 44294
$ws.Range("L433").Value = 'Primera'
$ws.Range("M433").Value = 80
$ws.Range("N433").Value = 7500
$ws.Range("O433").Value = 8000
$ws.Range("P433").Value = 7750
$ws.Range("R433").Value = 'Perú'
$ws.Range("S433").Value = 1938

# Row 434
$ws.Range("D434").Value = 44558
$ws.Range("L434").Value = 'Primera'
$ws.Range("M434").Value = 200
$ws.Range("N434").Value = 8500
$ws.Range("O434").Value = 9000
$ws.Range("P434").Value = 8750
$ws.Range("R434").Value = 'Perú'
$ws.Range("S434").Value = 2188

# Row 435
$ws.Range("D435").Value = 44558
$ws.Range("L435").Value = 'Segunda'
$ws.Range("M435").Value = 60
$ws.Range("N435").Value = 6000
$ws.Range("O435").Value = 6000
$ws.Range("P435").Value = 6000
$ws.Range("R435").Value = 'Perú'
$ws.Range("S435").Value = 1500

# Row 436
$ws.Range("D436").Value = 44775
$ws.Range("L436").Value = 'Primera'
$ws.Range("M436").Value = 200
$ws.Range("N436").Value = 13000
$ws.Range("O436").Value = 14000
$ws.Range("P436").Value = 13500
$ws.Range("R436").Value = 'Brasil'
$ws.Range("S436").Value = 3375

# Row 437
$ws.Range("D437").Value = 44859
$ws.Range("L437").Value = 'Primera'
$ws.Range("M437").Value = 100
$ws.Range("N437").Value = 9500
$ws.Range("O437").Value = 10000
$ws.Range("P437").Value = 9750
$ws.Range("R437").Value = 'Brasil'
$ws.Range("S437").Value = 2438

# Row 438
$ws.Range("D438").Value = 44859
$ws.Range("L438").Value = 'Segunda'
$ws.Range("M438").Value = 60
$ws.Range("N438").Value = 8000
$ws.Range("O438").Value = 8000
$ws.Range("P438").Value = 8000
$ws.Range("R438").Value = 'Brasil'
$ws.Range("S438").Value = 2000

# Row 439
$ws.Range("D439").Value = 44910
$ws.Range("L439").Value = 'Primera'
$ws.Range("M439").Value = 100
$ws.Range("N439").Value = 7500
$ws.Range("O439").Value = 8000
$ws.Range("P439").Value = 7750
$ws.Range("R439").Value = 'Brasil'
$ws.Range("S439").Value = 1938

# Row 440
$ws.Range("D440").Value = 45063
$ws.Range("L440").Value = 'Primera'
$ws.Range("M440").Value = 40
$ws.Range("N440").Value = 8000
$ws.Range("O440").Value = 8500
$ws.Range("P440").Value = 8250
$ws.Range("R440").Value = 'Perú'
$ws.Range("S440").Value = 2062

# Row 441
$ws.Range("D441").Value = 44365
$ws.Range("L441").Value = 'Primera'
$ws.Range("M441").Value = 160
$ws.Range("N441").Value = 7500
$ws.Range("O441").Value = 8000
$ws.Range("P441").Value = 7750
$ws.Range("R441").Value = 'Perú'
$ws.Range("S441").Value = 1938

# Row 442
$ws.Range("D442").Value = 45069
$ws.Range("L442").Value = 'Primera'
$ws.Range("M442").Value = 200
$ws.Range("N442").Value = 8000
$ws.Range("O442").Value = 8000
$ws.Range("P442").Value = 8000
$ws.Range("R442").Value = 'Perú'
$ws.Range("S442").Value = 2000

# Row 443
$ws.Range("D443").Value = 44795
$ws.Range("L443").Value = 'Primera'
$ws.Range("M443").Value = 200
$ws.Range("N443").Value = 13000
$ws.Range("O443").Value = 14000
$ws.Range("P443").Value = 13500
$ws.Range("R443").Value = 'Brasil'
$ws.Range("S443").Value = 3375

# Row 444
$ws.Range("D444").Value = 45173
$ws.Range("L444").Value = 'Primera'
$ws.Range("M444").Value = 40
$ws.Range("N444").Value = 12000
$ws.Range("O444").Value = 12000
$ws.Range("P444").Value = 12000
$ws.Range("R444").Value = 'Brasil'
$ws.Range("S444").Value = 3000

# Row 445
$ws.Range("D445").Value = 44673
$ws.Range("L445").Value = 'Primera'
$ws.Range("M445").Value = 200
$ws.Range("N445").Value = 7500
$ws.Range("O445").Value = 8000
$ws.Range("P445").Value = 7750
$ws.Range("R445").Value = 'Perú'
$ws.Range("S445").Value = 1938

# Row 446
$ws.Range("D446").Value = 44818
$ws.Range("L446").Value = 'Primera'
$ws.Range("M446").Value = 40
$ws.Range("N446").Value = 10000
$ws.Range("O446").Value = 11000
$ws.Range("P446").Value = 10500
$ws.Range("R446").Value = 'Brasil'
$ws.Range("S446").Value = 2625

# Row 447
$ws.Range("D447").Value = 45140
$ws.Range("L447").Value = 'Primera'
$ws.Range("M447").Value = 40
$ws.Range("N447").Value = 10000
$ws.Range("O447").Value = 10000
$ws.Range("P447").Value = 10000
$ws.Range("R447").Value = 'Perú'
$ws.Range("S447").Value = 2500

# Row 448
$ws.Range("D448").Value = 44649
$ws.Range("L448").Value = 'Primera'
$ws.Range("M448").Value = 200
$ws.Range("N448").Value = 7500
$ws.Range("O448").Value = 8000
$ws.Range("P448").Value = 7750
$ws.Range("R448").Value = 'Perú'
$ws.Range("S448").Value = 1938

# Row 449
$ws.Range("D449").Value = 45216
$ws.Range("L449").Value = 'Primera'
$ws.Range("M449").Value = 200
$ws.Range("N449").Value = 12000
$ws.Range("O449").Value = 12000
$ws.Range("P449").Value = 12000
$ws.Range("R449").Value = 'Brasil'
$ws.Range("S449").Value = 3000

# Row 450
$ws.Range("D450").Value = 44980
$ws.Range("L450").Value = 'Primera'
$ws.Range("M450").Value = 120
$ws.Range("N450").Value = 8000
$ws.Range("O450").Value = 9000
$ws.Range("P450").Value = 8500
$ws.Range("R450").Value = 'Perú'
$ws.Range("S450").Value = 2125

# Row 451
$ws.Range("D451").Value = 45114
$ws.Range("L451").Value = 'Primera'
$ws.Range("M451").Value = 200
$ws.Range("N451").Value = 8500
$ws.Range("O451").Value = 8500
$ws.Range("P451").Value = 8500
$ws.Range("R451").Value = 'Brasil'
$ws.Range("S451").Value = 2125

# Row 452
$ws.Range("D452").Value = 45012
$ws.Range("L452").Value = 'Primera'
$ws.Range("M452").Value = 120
$ws.Range("N452").Value = 8000
$ws.Range("O452").Value = 8500
$ws.Range("P452").Value = 8250
$ws.Range("R452").Value = 'Perú'
$ws.Range("S452").Value = 2062

# Row 453
$ws.Range("D453").Value = 44613
$ws.Range("L453").Value = 'Primera'
$ws.Range("M453").Value = 80
$ws.Range("N453").Value = 7500
$ws.Range("O453").Value = 7500
$ws.Range("P453").Value = 7500
$ws.Range("R453").Value = 'Perú'
$ws.Range("S453").Value = 1875

# Row 454
$ws.Range("D454").Value = 44893
$ws.Range("L454").Value = 'Primera'
$ws.Range("M454").Value = 100
$ws.Range("N454").Value = 9000
$ws.Range("O454").Value = 10000
$ws.Range("P454").Value = 9500
$ws.Range("R454").Value = 'Brasil'
$ws.Range("S454").Value = 2375

# Row 455
$ws.Range("D455").Value = 44392
$ws.Range("L455").Value = 'Primera'
$ws.Range("M455").Value = 150
$ws.Range("N455").Value = 7000
$ws.Range("O455").Value = 7000
$ws.Range("P455").Value = 7000
$ws.Range("R455").Value = 'Perú'
$ws.Range("S455").Value = 1750

# Row 456
$ws.Range("D456").Value = 44565
$ws.Range("L456").Value = 'Primera'
$ws.Range("M456").Value = 160
$ws.Range("N456").Value = 8500
$ws.Range("O456").Value = 9000
$ws.Range("P456").Value = 8750
$ws.Range("R456").Value = 'Perú'
$ws.Range("S456").Value = 2188

# Row 457
$ws.Range("D457").Value = 44565
$ws.Range("L457").Value = 'Segunda'
$ws.Range("M457").Value = 100
$ws.Range("N457").Value = 6000
$ws.Range("O457").Value = 6000
$ws.Range("P457").Value = 6000
$ws.Range("R457").Value = 'Perú'
$ws.Range("S457").Value = 1500

# Row 458
$ws.Range("D458").Value = 44544
$ws.Range("L458").Value = 'Primera'
$ws.Range("M458").Value = 200
$ws.Range("N458").Value = 7500
$ws.Range("O458").Value = 8000
$ws.Range("P458").Value = 7750
$ws.Range("R458").Value = 'Perú'
$ws.Range("S458").Value = 1938

# Row 459
$ws.Range("D459").Value = 44544
$ws.Range("L459").Value = 'Segunda'
$ws.Range("M459").Value = 100
$ws.Range("N459").Value = 5000
$ws.Range("O459").Value = 5000
$ws.Range("P459").Value = 5000
$ws.Range("R459").Value = 'Perú'
$ws.Range("S459").Value = 1250

# Fill in the columns for new row 459 that are constant across this dataset
$ws.Range("A459").Value = 4
$ws.Range("B459").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C459").Value = 'Los Lagos'
$ws.Range("E459").Value = 10
$ws.Range("F459").Value = 'Fruta'
$ws.Range("G459").Value = 100108
$ws.Range("H459").Value = 'Tropicales y subtropicales'
$ws.Range("I459").Value = 100108002
$ws.Range("J459").Value = 'Mango'
$ws.Range("K459").Value = 'Sin especificar'
$ws.Range("Q459").Value = '$/bandeja 4 kilos'
$ws.Range("T459").Value = 4

# Apply the same date/time number format used by the rest of column D to the new row
$ws.Range("D459").NumberFormat = $dateFormat